$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-21 05:07:06"
$wsZhCn.Range("H2").Value = "2016-03-21 05:07:46"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-21 05:07:14"
$wsDeDe.Range("H2").Value = "2016-03-21 05:08:00"
